$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1560
$ws.Range("I2").Value = 4373
$ws.Range("J2").Value = 17864
$ws.Range("K2").Value = 103
$ws.Range("L2").Value = 5009
$ws.Range("M2").Value = 303
$ws.Range("N2").Value = 3065
$ws.Range("O2").Value = 13
$ws.Range("P2").Value = 82
$ws.Range("Q2").Value = 23
$ws.Range("R2").Value = 249
$ws.Range("S2").Value = 1905
$ws.Range("T2").Value = 3238
$ws.Range("U2").Value = 243
$ws.Range("V2").Value = 27902
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 27762
$ws.Range("Y2").Value = 38
$ws.Range("Z2").Value = 427
$ws.Range("AA2").Value = 187
